$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as plain text, preserving General-like appearance afterward
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") '28.501.57'
Set-TextValue $ws.Range("E2") '  +1.22%  '
Set-TextValue $ws.Range("D3") '1.573.05'
Set-TextValue $ws.Range("E3") '  -0.21%  '
Set-TextValue $ws.Range("D5") '212.29'
Set-TextValue $ws.Range("E5") '  -0.34%  '
Set-TextValue $ws.Range("E7") '  +0.23%  '
Set-TextValue $ws.Range("D8") '46.29'
Set-TextValue $ws.Range("E8") '  +6.23%  '
Set-TextValue $ws.Range("D9") '24.15'
Set-TextValue $ws.Range("E9") '  +3.16%  '
Set-TextValue $ws.Range("E10") '  -1.35%  '
Set-TextValue $ws.Range("E11") '  -1.27%  '
Set-TextValue $ws.Range("E12") '  -0.18%  '
Set-TextValue $ws.Range("D13") '1.798.36'
Set-TextValue $ws.Range("E13") '  -0.19%  '
Set-TextValue $ws.Range("D14") '1.572.18'
Set-TextValue $ws.Range("E14") '  -0.27%  '
Set-TextValue $ws.Range("E15") '  -0.57%  '
Set-TextValue $ws.Range("E16") '  -1.49%  '
Set-TextValue $ws.Range("D17") '28.511.76'
Set-TextValue $ws.Range("E17") '  +1.41%  '
Set-TextValue $ws.Range("D18") '62.21'
Set-TextValue $ws.Range("E18") '  -2.21%  '
Set-TextValue $ws.Range("D19") '227.59'
Set-TextValue $ws.Range("D20") '7.37'
Set-TextValue $ws.Range("E20") '  -1.11%  '
Set-TextValue $ws.Range("D21") '0.0₃0693'
Set-TextValue $ws.Range("E21") '  -1.82%  '
Set-TextValue $ws.Range("E22") '  +0.26%  '
Set-TextValue $ws.Range("E23") '  -5.36%  '
Set-TextValue $ws.Range("E24") '  -2.42%  '
Set-TextValue $ws.Range("E25") '  +3.91%  '
Set-TextValue $ws.Range("D26") '151.32'
Set-TextValue $ws.Range("E26") '  -0.79%  '
Set-TextValue $ws.Range("D27") '14.99'
Set-TextValue $ws.Range("E27") '  -1.67%  '
Set-TextValue $ws.Range("E28") '  -1.93%  '
Set-TextValue $ws.Range("E29") '  -2.23%  '
Set-TextValue $ws.Range("E30") '  +0.23%  '
Set-TextValue $ws.Range("E31") '  -2.34%  '
Set-TextValue $ws.Range("E32") '  -2.05%  '
Set-TextValue $ws.Range("E33") '  -0.33%  '
Set-TextValue $ws.Range("E34") '  +0.54%  '
Set-TextValue $ws.Range("D35") '1.392.86'
Set-TextValue $ws.Range("E35") '  -1.73%  '
Set-TextValue $ws.Range("E36") '  -2.63%  '
Set-TextValue $ws.Range("E37") '  -2.20%  '
Set-TextValue $ws.Range("E38") '  +1.59%  '
Set-TextValue $ws.Range("D39") '2.62'
Set-TextValue $ws.Range("E39") '  +5.42%  '
Set-TextValue $ws.Range("E40") '  -0.85%  '
Set-TextValue $ws.Range("D41") '0.532'
Set-TextValue $ws.Range("E41") '  -1.72%  '
Set-TextValue $ws.Range("E42") '  +0.28%  '
Set-TextValue $ws.Range("E43") '  -1.62%  '
Set-TextValue $ws.Range("D44") '5.61'
Set-TextValue $ws.Range("E44") '  -0.23%  '
Set-TextValue $ws.Range("E45") '  +2.18%  '
Set-TextValue $ws.Range("E46") '  +0.81%  '
Set-TextValue $ws.Range("D47") '63.03'
Set-TextValue $ws.Range("E47") '  -1.12%  '
Set-TextValue $ws.Range("D48") '1.711.55'
Set-TextValue $ws.Range("E48") '  -0.22%  '
Set-TextValue $ws.Range("D49") '86.06'
Set-TextValue $ws.Range("E49") '  -1.23%  '
Set-TextValue $ws.Range("D50") '0.0₆0103'
Set-TextValue $ws.Range("E50") '  -3.37%  '
Set-TextValue $ws.Range("E51") '  -1.50%  '
